$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

$ws.Range("A3").Value = "TUHH_MECHATRONICS"
$ws.Range("B3").Value = "Yes"

$ws.Range("A4").Value = "HANNOVER_INTER_MECHATRONICS"
$ws.Range("B4").Value = "Yes"

$ws.Range("A5").Value = "TU_DORTMUND_MANUFACTURING_TECH"
$ws.Range("B5").Value = "Yes"

$ws.Range("A6").Select()
